# edit.ps1 - applies the SOP worksheet changes described by the commit:
#   "modify sop, add support for section, add support for increment/decrement on while"
#
# 1. Insert a new "section" header row right after the column-header row.
# 2. Change the "target criteria" value from "<=5  Å" to "lte 5  Å" everywhere it appears
#    (note: the space between "5" and "Å" is actually a THIN SPACE, U+2009).
# 3. Insert "flow operation"/"+" and "flow magnitude"/"1" rows into the step-6 "while"
#    loop block (right after its "flow compared value" row), mirroring the step-10
#    "for" loop block that already carries those two keys.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$thin = [char]0x2009
$aring = [char]0x00C5
$oldCriteria = "<=5 " + $thin + $aring
$newCriteria = "lte 5 " + $thin + $aring

# --- 1. Insert new section row at row 2 (pushes everything else down by one) ---
$ws.Rows.Item(2).EntireRow.Insert()
$ws.Cells.Item(2, 1).Value = "-"
$ws.Cells.Item(2, 2).Value = "section"
$ws.Cells.Item(2, 3).Value = "Structure Preparation"

# --- 2. Replace "<=5  Å" target-criteria value with "lte 5  Å" throughout column C ---
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cellVal = $cell.Value2
    if ($cellVal -eq $oldCriteria) {
        $cell.Value = $newCriteria
    }
}

# --- 3. Insert flow operation / flow magnitude rows into the step-6 "while" block ---
# After step (1), the step-6 "flow compared value" row (originally row 31) is now row 32.
$stepNum = $ws.Cells.Item(32, 1).Value2

$ws.Rows.Item(33).EntireRow.Insert()
$ws.Cells.Item(33, 1).Value = $stepNum
$ws.Cells.Item(33, 2).Value = "flow operation"
$ws.Cells.Item(33, 3).Value = "+"

$ws.Rows.Item(34).EntireRow.Insert()
$ws.Cells.Item(34, 1).Value = $stepNum
$ws.Cells.Item(34, 2).Value = "flow magnitude"
$ws.Cells.Item(34, 3).Value = "'1"

Write-Output "edit complete"
